$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.714.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.494.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.61"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.495.34"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.088.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.730.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.607"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.635.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000119"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.90%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.54%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.61"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.166"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.33"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.491.12"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.90"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0892"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "170.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.77%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.895"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.81"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.39%  "
